$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B" = 0.9999969141016266
    "C" = 0.9990169585194482
    "D" = 0.9999915046366025
    "E" = 0.999997548728656
    "F" = 0.9999936550366929
    "G" = 0.00000288055067624043
    "H" = 0.0009176260715613002
    "I" = 0.00001297647948492287
    "J" = 0.0000003349728090916505
    "K" = 0.000006655726147007262
    "L" = 0.00009999963621106398
    "M" = 0.00169721851163615
    "N" = 0.9999753128130129
    "O" = 0.001769472617403162
    "P" = 67.51505815013246
    "Q" = 93.11145047236465
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
